$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Rln1"
$ws.Cells.Item(2, 3).Value = "Rxfp1"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 2.119667333333334
$ws.Cells.Item(2, 8).Value = 6.359002
$ws.Cells.Item(2, 9).Value = 0.6371329247828699
$ws.Cells.Item(2, 10).Value = 0.6371329247828699
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 0.02528233333333333
$ws.Cells.Item(2, 14).Value = 0.075847
$ws.Cells.Item(2, 15).Value = 0.2893587312729617
$ws.Cells.Item(2, 16).Value = 0.2893587312729617
$ws.Cells.Item(2, 17).Value = 0.05359013607711112
$ws.Cells.Item(2, 18).Value = 0.482311224694
$ws.Cells.Item(2, 19).Value = 0.1843599747674026
$ws.Cells.Item(2, 20).Value = 0.1843599747674026

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Rln1"
$ws.Cells.Item(3, 3).Value = "Rxfp1"
$ws.Cells.Item(3, 4).Value = "MuSCs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 2.119667333333334
$ws.Cells.Item(3, 8).Value = 6.359002
$ws.Cells.Item(3, 9).Value = 0.6371329247828699
$ws.Cells.Item(3, 10).Value = 0.6371329247828699
$ws.Cells.Item(3, 11).Value = 2
$ws.Cells.Item(3, 12).Value = 0.6666666666666666
$ws.Cells.Item(3, 13).Value = 0.06209133333333333
$ws.Cells.Item(3, 14).Value = 0.186274
$ws.Cells.Item(3, 15).Value = 0.7106412687270383
$ws.Cells.Item(3, 16).Value = 0.7106412687270383
$ws.Cells.Item(3, 17).Value = 0.1316129709497778
$ws.Cells.Item(3, 18).Value = 1.184516738548
$ws.Cells.Item(3, 19).Value = 0.4527729500154673
$ws.Cells.Item(3, 20).Value = 0.4527729500154673

# Row 4
$ws.Cells.Item(4, 1).Value = "FAPs"
$ws.Cells.Item(4, 2).Value = "Rln1"
$ws.Cells.Item(4, 3).Value = "Rxfp1"
$ws.Cells.Item(4, 4).Value = "ECs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 0.8424356666666667
$ws.Cells.Item(4, 8).Value = 2.527307
$ws.Cells.Item(4, 9).Value = 0.2532206312773955
$ws.Cells.Item(4, 10).Value = 0.2532206312773955
$ws.Cells.Item(4, 11).Value = 1
$ws.Cells.Item(4, 12).Value = 0.3333333333333333
$ws.Cells.Item(4, 13).Value = 0.02528233333333333
$ws.Cells.Item(4, 14).Value = 0.075847
$ws.Cells.Item(4, 15).Value = 0.2893587312729617
$ws.Cells.Item(4, 16).Value = 0.2893587312729617
$ws.Cells.Item(4, 17).Value = 0.02129873933655556
$ws.Cells.Item(4, 18).Value = 0.191688654029
$ws.Cells.Item(4, 19).Value = 0.07327160059856561
$ws.Cells.Item(4, 20).Value = 0.07327160059856561

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Rln1"
$ws.Cells.Item(5, 3).Value = "Rxfp1"
$ws.Cells.Item(5, 4).Value = "MuSCs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 0.8424356666666667
$ws.Cells.Item(5, 8).Value = 2.527307
$ws.Cells.Item(5, 9).Value = 0.2532206312773955
$ws.Cells.Item(5, 10).Value = 0.2532206312773955
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 0.6666666666666666
$ws.Cells.Item(5, 13).Value = 0.06209133333333333
$ws.Cells.Item(5, 14).Value = 0.186274
$ws.Cells.Item(5, 15).Value = 0.7106412687270383
$ws.Cells.Item(5, 16).Value = 0.7106412687270383
$ws.Cells.Item(5, 17).Value = 0.05230795379088889
$ws.Cells.Item(5, 18).Value = 0.470771584118
$ws.Cells.Item(5, 19).Value = 0.1799490306788299
$ws.Cells.Item(5, 20).Value = 0.1799490306788299

# Row 6
$ws.Cells.Item(6, 1).Value = "MuSCs"
$ws.Cells.Item(6, 2).Value = "Rln1"
$ws.Cells.Item(6, 3).Value = "Rxfp1"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 0.2981566666666667
$ws.Cells.Item(6, 8).Value = 0.89447
$ws.Cells.Item(6, 9).Value = 0.08962039754517039
$ws.Cells.Item(6, 10).Value = 0.08962039754517039
$ws.Cells.Item(6, 11).Value = 1
$ws.Cells.Item(6, 12).Value = 0.3333333333333333
$ws.Cells.Item(6, 13).Value = 0.02528233333333333
$ws.Cells.Item(6, 14).Value = 0.075847
$ws.Cells.Item(6, 15).Value = 0.2893587312729617
$ws.Cells.Item(6, 16).Value = 0.2893587312729617
$ws.Cells.Item(6, 17).Value = 0.007538096232222223
$ws.Cells.Item(6, 18).Value = 0.06784286609
$ws.Cells.Item(6, 19).Value = 0.02593244452984896
$ws.Cells.Item(6, 20).Value = 0.02593244452984896

# Row 7
$ws.Cells.Item(7, 1).Value = "MuSCs"
$ws.Cells.Item(7, 2).Value = "Rln1"
$ws.Cells.Item(7, 3).Value = "Rxfp1"
$ws.Cells.Item(7, 4).Value = "MuSCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 0.2981566666666667
$ws.Cells.Item(7, 8).Value = 0.89447
$ws.Cells.Item(7, 9).Value = 0.08962039754517039
$ws.Cells.Item(7, 10).Value = 0.08962039754517039
$ws.Cells.Item(7, 11).Value = 2
$ws.Cells.Item(7, 12).Value = 0.6666666666666666
$ws.Cells.Item(7, 13).Value = 0.06209133333333333
$ws.Cells.Item(7, 14).Value = 0.186274
$ws.Cells.Item(7, 15).Value = 0.7106412687270383
$ws.Cells.Item(7, 16).Value = 0.7106412687270383
$ws.Cells.Item(7, 17).Value = 0.01851294497555556
$ws.Cells.Item(7, 18).Value = 0.16661650478
$ws.Cells.Item(7, 19).Value = 0.06368795301532143
$ws.Cells.Item(7, 20).Value = 0.06368795301532143

# Row 8
$ws.Cells.Item(8, 1).Value = "Resolving-Mac"
$ws.Cells.Item(8, 2).Value = "Rln1"
$ws.Cells.Item(8, 3).Value = "Rxfp1"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 0.06662433333333333
$ws.Cells.Item(8, 8).Value = 0.199873
$ws.Cells.Item(8, 9).Value = 0.0200260463945642
$ws.Cells.Item(8, 10).Value = 0.0200260463945642
$ws.Cells.Item(8, 11).Value = 1
$ws.Cells.Item(8, 12).Value = 0.3333333333333333
$ws.Cells.Item(8, 13).Value = 0.02528233333333333
$ws.Cells.Item(8, 14).Value = 0.075847
$ws.Cells.Item(8, 15).Value = 0.2893587312729617
$ws.Cells.Item(8, 16).Value = 0.2893587312729617
$ws.Cells.Item(8, 17).Value = 0.001684418603444444
$ws.Cells.Item(8, 18).Value = 0.015159767431
$ws.Cells.Item(8, 19).Value = 0.005794711377144566
$ws.Cells.Item(8, 20).Value = 0.005794711377144567

# Row 9
$ws.Cells.Item(9, 1).Value = "Resolving-Mac"
$ws.Cells.Item(9, 2).Value = "Rln1"
$ws.Cells.Item(9, 3).Value = "Rxfp1"
$ws.Cells.Item(9, 4).Value = "MuSCs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 0.06662433333333333
$ws.Cells.Item(9, 8).Value = 0.199873
$ws.Cells.Item(9, 9).Value = 0.0200260463945642
$ws.Cells.Item(9, 10).Value = 0.0200260463945642
$ws.Cells.Item(9, 11).Value = 2
$ws.Cells.Item(9, 12).Value = 0.6666666666666666
$ws.Cells.Item(9, 13).Value = 0.06209133333333333
$ws.Cells.Item(9, 14).Value = 0.186274
$ws.Cells.Item(9, 15).Value = 0.7106412687270383
$ws.Cells.Item(9, 16).Value = 0.7106412687270383
$ws.Cells.Item(9, 17).Value = 0.004136793689111111
$ws.Cells.Item(9, 18).Value = 0.037231143202
$ws.Cells.Item(9, 19).Value = 0.01423133501741963
$ws.Cells.Item(9, 20).Value = 0.01423133501741963

